# "Generate Report for handback"
#
# For each locale sheet (zh-cn, de-de) the localization files listed in rows 2
# and 3 have now been handed back from translation. The report needs to:
#   - flip the Status column (B) from "Not yet handed off" to "Handed back"
#   - populate the "Latest Target File" (E) and "Latest Handback File" (F)
#     columns with the same source/handoff file references as columns A/C
#     (as hyperlinks, matching the existing hyperlink styling/behaviour)
#   - stamp the "Latest Handback DateTime" (G) with the handback timestamp

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = "Handed back"
$wsZh.Range("B3").Value = "Handed back"

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/a7688c13fa71b1bc5fa376033f50ede434add8c0/e2e/4694cac1-0c76-460a-9e37-eed913adbefa.md", "", "", "4694cac1-0c76-460a-9e37-eed913adbefa.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dbad3b68175a3135e0c7dc4b67136c9b43ec7e07/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/4694cac1-0c76-460a-9e37-eed913adbefa.e3e64a72f20027cc0ee51f5dd4322b38e3cf3023.zh-cn.xlf", "", "", "4694cac1-0c76-460a-9e37-eed913adbefa.e3e64a72f20027cc0ee51f5dd4322b38e3cf3023.zh-cn.xlf")
$wsZh.Range("G2").Value = "2016-01-08 09:02:15"

$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/a7688c13fa71b1bc5fa376033f50ede434add8c0/e2e/5c8a57b5-dcc1-4db6-ab2e-fb1263056230.md", "", "", "5c8a57b5-dcc1-4db6-ab2e-fb1263056230.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dbad3b68175a3135e0c7dc4b67136c9b43ec7e07/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/5c8a57b5-dcc1-4db6-ab2e-fb1263056230.b89e86b7a1638172d03105f9538774ae505a43c2.zh-cn.xlf", "", "", "5c8a57b5-dcc1-4db6-ab2e-fb1263056230.b89e86b7a1638172d03105f9538774ae505a43c2.zh-cn.xlf")
$wsZh.Range("G3").Value = "2016-01-08 09:02:15"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = "Handed back"
$wsDe.Range("B3").Value = "Handed back"

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/a7688c13fa71b1bc5fa376033f50ede434add8c0/e2e/4694cac1-0c76-460a-9e37-eed913adbefa.md", "", "", "4694cac1-0c76-460a-9e37-eed913adbefa.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6daadcf9f95d886d615499a9f8e9d2f3451fa4ab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/4694cac1-0c76-460a-9e37-eed913adbefa.e3e64a72f20027cc0ee51f5dd4322b38e3cf3023.de-de.xlf", "", "", "4694cac1-0c76-460a-9e37-eed913adbefa.e3e64a72f20027cc0ee51f5dd4322b38e3cf3023.de-de.xlf")
$wsDe.Range("G2").Value = "2016-01-08 09:02:32"

$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/a7688c13fa71b1bc5fa376033f50ede434add8c0/e2e/5c8a57b5-dcc1-4db6-ab2e-fb1263056230.md", "", "", "5c8a57b5-dcc1-4db6-ab2e-fb1263056230.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6daadcf9f95d886d615499a9f8e9d2f3451fa4ab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/5c8a57b5-dcc1-4db6-ab2e-fb1263056230.b89e86b7a1638172d03105f9538774ae505a43c2.de-de.xlf", "", "", "5c8a57b5-dcc1-4db6-ab2e-fb1263056230.b89e86b7a1638172d03105f9538774ae505a43c2.de-de.xlf")
$wsDe.Range("G3").Value = "2016-01-08 09:02:32"
